$wb = $excel.ActiveWorkbook

# --- Rename sheets: TP -> LL, FP -> NL ---
$wsLL = $wb.Worksheets.Item("TP")
$wsLL.Name = "LL"

$wsNL = $wb.Worksheets.Item("FP")
$wsNL.Name = "NL"

$wsTRUTH = $wb.Worksheets.Item("TRUTH")

# --- Rename rating columns to the unified "NLRating" header ---
$wsLL.Range("E1").Value = "NLRating"
$wsNL.Range("D1").Value = "NLRating"

# --- Drop the Paradigm/FROC/FCTRL columns (D:F) from the TRUTH sheet ---
$wsTRUTH.Columns("D:F").Delete()

# --- Restore view/selection state: LL sheet selection, TRUTH sheet selection, NL active with its own selection ---
$wsLL.Select()
$wsLL.Range("E2").Select()

$wsTRUTH.Select()
$wsTRUTH.Range("D1:F1048576").Select()

$wsNL.Select()
$wsNL.Range("D2").Select()
